$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the DATE for the existing (previously blank) second row
$ws.Range("A2").Value = 45928
$ws.Range("A2").NumberFormat = "mm-dd-yy"

# Grow Table1 by one row (adds a new formatted row at the bottom of the table)
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

# Copy row 2's formatting down into the freshly-added row 3, then fill in its values
$ws.Range("A2:C2").Copy($ws.Range("A3:C3"))
$ws.Range("A3").Value = 45934
$ws.Range("B3").Value = "Be Thou My Vision, Tsis Muaj Koj Pab (#148), You Raise Me Up, 10,000 Reasons"
$ws.Range("C3").Value = "40th Year Anniversary"

# Column widths widened (best-fit) to accommodate the new, longer content
$ws.Columns.Item(1).ColumnWidth = 8.833333333333332
$ws.Columns.Item(2).ColumnWidth = 70.5
$ws.Columns.Item(3).ColumnWidth = 19.666666666666668

# Matches the saved selection in the target workbook
$ws.Range("B8").Select() | Out-Null
